$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (2020 cohort, period 6): num_customers 7 -> 8, retention_rate recalculated
$ws.Range("C16").Value = 8
$ws.Range("E16").Value = 0.003796867584243

# Row 27 (2021 cohort, period 4): num_customers 34 -> 35, retention_rate recalculated
$ws.Range("C27").Value = 35
$ws.Range("E27").Value = 0.01554174067495559

# Row 31 (2022 cohort, period 3): num_customers 34 -> 35, retention_rate recalculated
$ws.Range("C31").Value = 35
$ws.Range("E31").Value = 0.01513840830449827

# Row 34 (2023 cohort, period 2): num_customers 53 -> 54, retention_rate recalculated
$ws.Range("C34").Value = 54
$ws.Range("E34").Value = 0.02393617021276596

# Row 36 (2024 cohort, period 1): num_customers 95 -> 97, retention_rate recalculated
$ws.Range("C36").Value = 97
$ws.Range("E36").Value = 0.05025906735751295

# Row 37 (2025 cohort, period 0): num_customers 514 -> 552, cohort_size 514 -> 552
$ws.Range("C37").Value = 552
$ws.Range("D37").Value = 552
